$d = $word.ActiveDocument

$replacements = @(
    @("399÷5=", "261÷7="),
    @("881÷2=", "469÷3="),
    @("451÷7=", "794÷6="),
    @("310÷9=", "514÷2="),
    @("354÷4=", "771÷6="),
    @("255÷2=", "514÷7="),
    @("931÷6=", "971÷7="),
    @("324÷7=", "196÷3="),
    @("791÷7=", "357÷8="),
    @("189÷6=", "591÷2="),
    @("564÷6=", "731÷3="),
    @("117÷8=", "283÷6="),
    @("524÷6=", "404÷7="),
    @("539÷5=", "106÷6="),
    @("902÷4=", "180÷5="),
    @("234÷6=", "104÷4="),
    @("562÷4=", "503÷6="),
    @("215÷4=", "368÷8="),
    @("323÷6=", "275÷9="),
    @("869÷5=", "476÷9="),
    @("288÷3=", "536÷4="),
    @("217÷4=", "720÷9="),
    @("292÷9=", "840÷8="),
    @("538÷9=", "792÷6="),
    @("466÷6=", "578÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
